$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 64009.258
$ws.Range("J17").Value = 65877.17999999999
$ws.Range("L17").Value = 197631.54
$ws.Range("N17").Value = -197967.54

$ws.Range("H57").Value = 35693.332
$ws.Range("J57").Value = 33640
$ws.Range("L57").Value = 100920
$ws.Range("N57").Value = -101918

$ws.Range("H107").Value = 909.56525
$ws.Range("I107").Value = 928
$ws.Range("K107").Value = 928
$ws.Range("M107").Value = 992

$ws.Range("H113").Value = 3668.375
$ws.Range("I113").Value = 3950
$ws.Range("J113").Value = 3449.3333
$ws.Range("K113").Value = 3950
$ws.Range("L113").Value = 3449.3333
$ws.Range("M113").Value = -696
$ws.Range("N113").Value = -9957.3333

$ws.Range("H132").Value = 6063764
$ws.Range("I132").Value = 6454632.5
$ws.Range("K132").Value = 19363897.5
$ws.Range("M132").Value = -19361367.5

$ws.Range("H137").Value = 3842
$ws.Range("I137").Value = 3973.1785
$ws.Range("J137").Value = 3535.9167
$ws.Range("K137").Value = 11919.5355
$ws.Range("L137").Value = 10607.7501
$ws.Range("M137").Value = -9369.5355
$ws.Range("N137").Value = -15707.7501

$ws.Range("H138").Value = 3662.6187
$ws.Range("I138").Value = 1757.4117
$ws.Range("J138").Value = 4690.825
$ws.Range("K138").Value = 5272.2351
$ws.Range("L138").Value = 14072.475
$ws.Range("M138").Value = -132.2350999999999
$ws.Range("N138").Value = -24352.475

$ws.Range("H140").Value = 29450
$ws.Range("J140").Value = 29450
$ws.Range("L140").Value = 29450
$ws.Range("N140").Value = -39810

$ws.Range("H141").Value = 405145.16
$ws.Range("I141").Value = 1555.3334
$ws.Range("J141").Value = 539675.1
$ws.Range("K141").Value = 4666.0002
$ws.Range("L141").Value = 1619025.3
$ws.Range("M141").Value = 513.9997999999996
$ws.Range("N141").Value = -1629385.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 14704.9
$ws.Range("J33").Value = 16289.857
$ws.Range("L33").Value = 16289.857
$ws.Range("N33").Value = -16947.857

$ws.Range("H74").Value = 1482.7963
$ws.Range("I74").Value = 1181.8445
$ws.Range("J74").Value = 2987.5557
$ws.Range("K74").Value = 1181.8445
$ws.Range("L74").Value = 2987.5557
$ws.Range("M74").Value = -307.8444999999999
$ws.Range("N74").Value = -4735.5557

$ws.Range("H77").Value = 1482.7963
$ws.Range("I77").Value = 1181.8445
$ws.Range("J77").Value = 2987.5557
$ws.Range("K77").Value = 5909.2225
$ws.Range("L77").Value = 14937.7785
$ws.Range("M77").Value = -1541.2225
$ws.Range("N77").Value = -23673.7785

$ws.Range("H110").Value = 1449.75
$ws.Range("I110").Value = 560.8077
$ws.Range("K110").Value = 560.8077
$ws.Range("M110").Value = 1484.1923

$ws.Range("H134").Value = 30290
$ws.Range("J134").Value = 30290
$ws.Range("L134").Value = 30290
$ws.Range("N134").Value = -40430

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 68423.664
$ws.Range("I23").Value = 25250
$ws.Range("K23").Value = 25250
$ws.Range("M23").Value = -24967

$ws.Range("H74").Value = 51745
$ws.Range("J74").Value = 51745
$ws.Range("L74").Value = 51745
$ws.Range("N74").Value = -53617

$ws.Range("H77").Value = 51745
$ws.Range("J77").Value = 51745
$ws.Range("L77").Value = 155235
$ws.Range("N77").Value = -164595

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4059.5667
$ws.Range("I62").Value = 2376.5
$ws.Range("J62").Value = 6584.1665
$ws.Range("K62").Value = 2376.5
$ws.Range("L62").Value = 6584.1665
$ws.Range("M62").Value = -1752.5
$ws.Range("N62").Value = -7832.1665

$ws.Range("H65").Value = 4059.5667
$ws.Range("I65").Value = 2376.5
$ws.Range("J65").Value = 6584.1665
$ws.Range("K65").Value = 11882.5
$ws.Range("L65").Value = 32920.8325
$ws.Range("M65").Value = -8762.5
$ws.Range("N65").Value = -39160.8325

$ws.Range("H107").Value = 1229.1111
$ws.Range("I107").Value = 1007.3889
$ws.Range("K107").Value = 1007.3889
$ws.Range("M107").Value = 912.6111

$ws.Range("H132").Value = 1911.826
$ws.Range("I132").Value = 1634.6111
$ws.Range("J132").Value = 2909.8
$ws.Range("K132").Value = 4903.8333
$ws.Range("L132").Value = 8729.400000000001
$ws.Range("M132").Value = -2373.8333
$ws.Range("N132").Value = -13789.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 30550.723
$ws.Range("I129").Value = 3907.1428
$ws.Range("J129").Value = 47505.727
$ws.Range("K129").Value = 11721.4284
$ws.Range("L129").Value = 142517.181
$ws.Range("M129").Value = -6721.428400000001
$ws.Range("N129").Value = -152517.181

$ws.Range("H131").Value = 1360.4445
$ws.Range("I131").Value = 2606
$ws.Range("J131").Value = 1125.434
$ws.Range("K131").Value = 7818
$ws.Range("L131").Value = 3376.302
$ws.Range("M131").Value = -2778
$ws.Range("N131").Value = -13456.302

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3459.5557
$ws.Range("I132").Value = 3454.258
$ws.Range("J132").Value = 3471.2856
$ws.Range("K132").Value = 10362.774
$ws.Range("L132").Value = 10413.8568
$ws.Range("M132").Value = -7832.773999999999
$ws.Range("N132").Value = -15473.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 166673660
$ws.Range("J61").Value = 10666.667
$ws.Range("L61").Value = 10666.667
$ws.Range("N61").Value = -11070.667

$ws.Range("H113").Value = 166673660
$ws.Range("J113").Value = 10666.667
$ws.Range("L113").Value = 10666.667
$ws.Range("N113").Value = -15006.667

$ws.Range("H132").Value = 2413.9148
$ws.Range("I132").Value = 1691.6666
$ws.Range("J132").Value = 3688.4707
$ws.Range("K132").Value = 5074.9998
$ws.Range("L132").Value = 11065.4121
$ws.Range("M132").Value = -2544.9998
$ws.Range("N132").Value = -16125.4121

$ws.Range("H136").Value = 4507.6787
$ws.Range("I136").Value = 3713.125
$ws.Range("K136").Value = 11139.375
$ws.Range("M136").Value = -8589.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2000742.4
$ws.Range("I3").Value = 10000000
$ws.Range("K3").Value = 10000000
$ws.Range("M3").Value = -9999886

$ws.Range("H113").Value = 894.6429000000001
$ws.Range("I113").Value = 176
$ws.Range("J113").Value = 1613.2858
$ws.Range("K113").Value = 528
$ws.Range("L113").Value = 4839.857400000001
$ws.Range("M113").Value = 1642
$ws.Range("N113").Value = -9179.857400000001
